# ---------------------------------------------------------------------------
# Edit summary (per the target OOXML diff):
#
# 1. Slide 5 contains a single table (graphicFrame, Shapes.Item(2)); its
#    table style is switched from the custom "Table_0" style
#    {9A4DB9A1-787B-42AC-BB83-BE1123A7C42C} to the built-in style
#    {B541A4C5-26D9-4012-A64D-E28A12D294B9}. Table styles can't be assigned
#    through the Style property directly (PowerPoint requires ApplyStyle).
#
# 2. The deck's active theme ("Integral" / "Red Violet" colour scheme, used
#    by the slide master/theme2.xml) is swapped out for the classic
#    "Office" colour scheme (the colours that originally lived in the
#    otherwise-unused theme1.xml, which is only wired to the notes master).
#    We reproduce that swap at the level that is actually exposed through
#    the PowerPoint object model: the twelve theme colours themselves,
#    via Slide.ThemeColorScheme (which edits the shared theme part used by
#    the slide master, i.e. theme2.xml).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style -----------------------------------------------------

$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{B541A4C5-26D9-4012-A64D-E28A12D294B9}")

# --- 2. Theme colours -----------------------------------------------------
# Order exposed by ThemeColorScheme.Item(n) is fixed:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# RGB() builds the BGR-packed long that PowerPoint's ColorFormat.RGB expects.

function RGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$officeColors = @(
    (RGB 0x00 0x00 0x00),  # dk1      000000
    (RGB 0xFF 0xFF 0xFF),  # lt1      FFFFFF
    (RGB 0x44 0x54 0x6A),  # dk2      44546A
    (RGB 0xE7 0xE6 0xE6),  # lt2      E7E6E6
    (RGB 0x5B 0x9B 0xD5),  # accent1  5B9BD5
    (RGB 0xED 0x7D 0x31),  # accent2  ED7D31
    (RGB 0xA5 0xA5 0xA5),  # accent3  A5A5A5
    (RGB 0xFF 0xC0 0x00),  # accent4  FFC000
    (RGB 0x44 0x72 0xC4),  # accent5  4472C4
    (RGB 0x70 0xAD 0x47),  # accent6  70AD47
    (RGB 0x05 0x63 0xC1),  # hlink    0563C1
    (RGB 0x95 0x4F 0x72)   # folHlink 954F72
)

$themeSlide = $p.Slides.Item(1)
$colorScheme = $themeSlide.ThemeColorScheme

for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Item($i).RGB = $officeColors[$i - 1]
}
